$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.088045
$ws.Range("H2").Value = 0.264135
$ws.Range("I2").Value = 0.1003389266487061
$ws.Range("J2").Value = 0.1003389266487061
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 0.2398347267416667
$ws.Range("R2").Value = 2.158512540675
$ws.Range("S2").Value = 0.004651211420644677
$ws.Range("T2").Value = 0.004651211420644677
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.088045
$ws.Range("H3").Value = 0.264135
$ws.Range("I3").Value = 0.1003389266487061
$ws.Range("J3").Value = 0.1003389266487061
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 3.576443427046667
$ws.Range("R3").Value = 32.18799084342
$ws.Range("S3").Value = 0.06935940736842039
$ws.Range("T3").Value = 0.06935940736842039
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.088045
$ws.Range("H4").Value = 0.264135
$ws.Range("I4").Value = 0.1003389266487061
$ws.Range("J4").Value = 0.1003389266487061
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 1.357590947825
$ws.Range("R4").Value = 12.218318530425
$ws.Range("S4").Value = 0.02632830785964099
$ws.Range("T4").Value = 0.02632830785964099
$ws.Range("G5").Value = 0.3888126666666667
$ws.Range("I5").Value = 0.4431034770941504
$ws.Range("J5").Value = 0.4431034770941504
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 1.059126352021111
$ws.Range("R5").Value = 9.532137168190001
$ws.Range("S5").Value = 0.02054006378205818
$ws.Range("T5").Value = 0.02054006378205818
$ws.Range("G6").Value = 0.3888126666666667
$ws.Range("I6").Value = 0.4431034770941504
$ws.Range("J6").Value = 0.4431034770941504
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.3062958275578986
$ws.Range("T6").Value = 0.3062958275578986
$ws.Range("G7").Value = 0.3888126666666667
$ws.Range("I7").Value = 0.4431034770941504
$ws.Range("J7").Value = 0.4431034770941504
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 5.995213318943334
$ws.Range("R7").Value = 53.95691987049
$ws.Range("S7").Value = 0.1162675857541936
$ws.Range("T7").Value = 0.1162675857541936
$ws.Range("G8").Value = 0.4006183333333334
$ws.Range("H8").Value = 1.201855
$ws.Range("I8").Value = 0.4565575962571436
$ws.Range("J8").Value = 0.4565575962571436
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 1.091285007697222
$ws.Range("R8").Value = 9.821565069275001
$ws.Range("S8").Value = 0.02116372953966308
$ws.Range("T8").Value = 0.02116372953966308
$ws.Range("G9").Value = 0.4006183333333334
$ws.Range("H9").Value = 1.201855
$ws.Range("I9").Value = 0.4565575962571436
$ws.Range("J9").Value = 0.4565575962571436
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 16.27336935662889
$ws.Range("R9").Value = 146.46032420966
$ws.Range("S9").Value = 0.3155960040993162
$ws.Range("T9").Value = 0.3155960040993162
$ws.Range("G10").Value = 0.4006183333333334
$ws.Range("H10").Value = 1.201855
$ws.Range("I10").Value = 0.4565575962571436
$ws.Range("J10").Value = 0.4565575962571436
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 6.177248257891667
$ws.Range("R10").Value = 55.595234321025
$ws.Range("S10").Value = 0.1197978626181643
$ws.Range("T10").Value = 0.1197978626181643
